$wb = $excel.ActiveWorkbook

# --- Update the "Structures" sheet's System column (B) to use full system
# --- names instead of the short keyword, for all 4 data rows.
$wsStructures = $wb.Worksheets.Item("Structures")
$wsStructures.Range("B2").Value = "Nintendo Entertainment System"
$wsStructures.Range("B3").Value = "Super Nintendo Entertainment System"
$wsStructures.Range("B4").Value = "Nintendo Entertainment System"
$wsStructures.Range("B5").Value = "Super Nintendo Entertainment System"

# --- Update selections/active cells on both sheets.
$wsSystems = $wb.Worksheets.Item("Systems")
$wsSystems.Range("A4").Select() | Out-Null

$wsStructures.Activate() | Out-Null
$wsStructures.Range("B5").Select() | Out-Null
